# Update the NATMI LR-pair (Efnb2-Epha3) worksheet with newly computed TPM-based
# statistics. Only the numeric result columns (G, H, I, J, M, N, O, P, Q, R, S, T)
# change for the relevant rows; identifiers/categorical columns are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 42.75280866666666
$ws.Range("H2").Value = 128.258426
$ws.Range("I2").Value = 0.8529286054750734
$ws.Range("J2").Value = 0.8529286054750735
$ws.Range("M2").Value = 0.003058333333333333
$ws.Range("N2").Value = 0.009175
$ws.Range("O2").Value = 0.0001379486413073712
$ws.Range("P2").Value = 0.0001379486413073712
$ws.Range("Q2").Value = 0.1307523398388889
$ws.Range("R2").Value = 1.17677105855
$ws.Range("S2").Value = 0.0001176603422574773
$ws.Range("T2").Value = 0.0001176603422574773
$ws.Range("G3").Value = 42.75280866666666
$ws.Range("H3").Value = 128.258426
$ws.Range("I3").Value = 0.8529286054750734
$ws.Range("J3").Value = 0.8529286054750735
$ws.Range("O3").Value = 0.9939610820947024
$ws.Range("P3").Value = 0.9939610820947024
$ws.Range("Q3").Value = 942.1095848497621
$ws.Range("R3").Value = 8478.98626364786
$ws.Range("S3").Value = 0.8477778396475295
$ws.Range("T3").Value = 0.8477778396475296
$ws.Range("G4").Value = 42.75280866666666
$ws.Range("H4").Value = 128.258426
$ws.Range("I4").Value = 0.8529286054750734
$ws.Range("J4").Value = 0.8529286054750735
$ws.Range("O4").Value = 0.005900969263990248
$ws.Range("P4").Value = 0.005900969263990248
$ws.Range("Q4").Value = 5.593136193816665
$ws.Range("R4").Value = 50.33822574435
$ws.Range("S4").Value = 0.005033105485286473
$ws.Range("T4").Value = 0.005033105485286473
$ws.Range("I5").Value = 0.04642608686423023
$ws.Range("J5").Value = 0.04642608686423023
$ws.Range("M5").Value = 0.003058333333333333
$ws.Range("N5").Value = 0.009175
$ws.Range("O5").Value = 0.0001379486413073712
$ws.Range("P5").Value = 0.0001379486413073712
$ws.Range("Q5").Value = 0.00711703118888889
$ws.Range("R5").Value = 0.0640532807
$ws.Range("S5").Value = 0.000006404415604138555
$ws.Range("T5").Value = 0.000006404415604138555
$ws.Range("I6").Value = 0.04642608686423023
$ws.Range("J6").Value = 0.04642608686423023
$ws.Range("O6").Value = 0.9939610820947024
$ws.Range("P6").Value = 0.9939610820947024
$ws.Range("S6").Value = 0.04614572353699292
$ws.Range("T6").Value = 0.04614572353699292
$ws.Range("I7").Value = 0.04642608686423023
$ws.Range("J7").Value = 0.04642608686423023
$ws.Range("O7").Value = 0.005900969263990248
$ws.Range("P7").Value = 0.005900969263990248
$ws.Range("S7").Value = 0.000273958911633164
$ws.Range("T7").Value = 0.000273958911633164
$ws.Range("G8").Value = 5.044817999999999
$ws.Range("I8").Value = 0.1006453076606963
$ws.Range("J8").Value = 0.1006453076606963
$ws.Range("M8").Value = 0.003058333333333333
$ws.Range("N8").Value = 0.009175
$ws.Range("O8").Value = 0.0001379486413073712
$ws.Range("P8").Value = 0.0001379486413073712
$ws.Range("Q8").Value = 0.01542873505
$ws.Range("R8").Value = 0.13885861545
$ws.Range("S8").Value = 0.00001388388344575541
$ws.Range("T8").Value = 0.00001388388344575541
$ws.Range("G9").Value = 5.044817999999999
$ws.Range("I9").Value = 0.1006453076606963
$ws.Range("J9").Value = 0.1006453076606963
$ws.Range("O9").Value = 0.9939610820947024
$ws.Range("P9").Value = 0.9939610820947024
$ws.Range("S9").Value = 0.1000375189101799
$ws.Range("T9").Value = 0.1000375189101799
$ws.Range("G10").Value = 5.044817999999999
$ws.Range("I10").Value = 0.1006453076606963
$ws.Range("J10").Value = 0.1006453076606963
$ws.Range("O10").Value = 0.005900969263990248
$ws.Range("P10").Value = 0.005900969263990248
$ws.Range("Q10").Value = 0.6599883148499999
$ws.Range("R10").Value = 5.939894833649999
$ws.Range("S10").Value = 0.0005939048670706111
$ws.Range("T10").Value = 0.0005939048670706111
